# ManageProducts.xlsx - "classic view check addition and order date
# selection code update"
#
# The QA fixture's "Input" sheet keeps one randomly generated SKU value
# per pallet-type row (B2:B7). A fresh batch of SKU values was rolled for
# the NP-SC-SKU (B2), NP-MC-SKU (B3) and P-MC-SKU (B5) rows, and the
# "classic view" cell styling (thin top+bottom border, solid white
# fill) was re-stamped on each of those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ClassicViewCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.Value = $value

    # "classic view" look: thin border on top + bottom, solid white fill
    $rng.Borders(8).LineStyle = 1   # xlEdgeTop
    $rng.Borders(8).Weight = 2      # xlThin
    $rng.Borders(9).LineStyle = 1   # xlEdgeBottom
    $rng.Borders(9).Weight = 2      # xlThin
    $rng.Interior.Pattern = 1       # xlSolid
    $rng.Interior.Color = 16777215  # RGB(255,255,255) -> indexed white
}

# New order-date-selection SKUs, re-rolled in row order.
Set-ClassicViewCell "B3" "prodClbw"
Set-ClassicViewCell "B5" "prodFeJw"
Set-ClassicViewCell "B2" "prodEsCE"
